# Insert a new data row at row 32 (pushes existing rows 32..152 down to 33..153)
# and populate it with the new record described by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(32).Insert()

$ws.Cells.Item(32, 1).Value = 10
$ws.Cells.Item(32, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(32, 3).Value = "La Araucanía"
$ws.Cells.Item(32, 4).Value = 44525
$ws.Cells.Item(32, 5).Value = 9
$ws.Cells.Item(32, 6).Value = 100112052
$ws.Cells.Item(32, 7).Value = "Albahaca"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 40
$ws.Cells.Item(32, 11).Value = 8000
$ws.Cells.Item(32, 12).Value = 8000
$ws.Cells.Item(32, 13).Value = 8000
$ws.Cells.Item(32, 14).Value = "`$/paquete"
$ws.Cells.Item(32, 15).Value = "Región del Maule"
$ws.Cells.Item(32, 16).Value = 8000
$ws.Cells.Item(32, 17).Value = 1
$ws.Cells.Item(32, 18).Value = "Hortaliza"
